{"js": "// Locate the first paragraph (\"This is a Microsoft word document.\") and\n// append \" (Changed main)\" to it as three additional runs, matching the\n// target OOXML diff:\n//   <w:r><w:t>This is a Microsoft word document.</w:t></w:r>\n//   <w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>\n//   <w:r><w:t>Changed main</w:t></w:r>\n//   <w:r><w:t>)</w:t></w:r>\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst target = paragraphs.items[0];\ntarget.load(\"text\");\nawait context.sync();\n\n// Sanity check - only touch the intended paragraph.\nif (target.text.indexOf(\"This is a Microsoft word document.\") !== 0) {\n  throw new Error(\"Expected first paragraph to start with the known sentence.\");\n}\n\n// Insert three separate runs at the end of the paragraph via OOXML so each\n// piece of text ends up as its own <w:r>, exactly like the authored edit.\nconst endRange = target.getRange(\"End\");\n\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>' +\n  '<w:r><w:t>Changed main</w:t></w:r>' +\n  '<w:r><w:t>)</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nendRange.insertOoxml(ooxml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Locate the first paragraph's text (\"This is a Microsoft word document.\")\n# and insert \" (Changed main)\" right after it as three additional runs,\n# matching the target OOXML diff:\n#   <w:r><w:t>This is a Microsoft word document.</w:t></w:r>\n#   <w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>\n#   <w:r><w:t>Changed main</w:t></w:r>\n#   <w:r><w:t>)</w:t></w:r>\n\n$d = $word.ActiveDocument\n\n$target = $d.Content\n$find = $target.Find\n$find.Text = \"This is a Microsoft word document.\"\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find the target sentence 'This is a Microsoft word document.'\"\n}\n\n# $target now spans exactly the matched sentence (no trailing paragraph\n# mark). Feeding Range.InsertXML a WordprocessingML package whose body\n# starts with the same sentence (re-asserted as its own run) followed by\n# the three new runs replaces the range with four distinct <w:r> elements,\n# rather than merging the new text into the existing run.\n$ooxml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t>This is a Microsoft word document.</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>\n            <w:r><w:t>Changed main</w:t></w:r>\n            <w:r><w:t>)</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$target.InsertXML($ooxml)\n"}
